# "Generate Report for Handback" - refresh the localization-status report:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     on the Overview sheet and on each per-language detail sheet.
#   - The per-language "Latest Handback DateTime" is refreshed to the new
#     handback timestamp.
#   - The stale "handback file is not the latest" Error Detail is cleared now
#     that the handback is in sync.
#   - Column widths are widened/narrowed to fit the new Status text / the now
#     (mostly) empty Error Detail column.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# --- zh-cn detail sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("K2").Value = "2016-09-03 10:53:57"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(16).ColumnWidth = 12.75

# --- de-de detail sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("K2").Value = "2016-09-03 10:54:10"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(16).ColumnWidth = 12.75
